# TC38_Canine_Filter_Breed-Poodle.xlsx — "updated the queries to fix case
# files comparison icdc"
#
# The FilesTab Cypher query stored in B4 is rewritten: the column order of
# the RETURN clause changes (Format/File Type/Size move earlier, Association
# moves later) and a stray blank line right after the WHERE clause is
# removed. The author then left the cursor/selection sitting on B4 (instead
# of the previous B2) when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFileQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle'] 
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
       coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

# Rewrite the FilesTab query text in place.
$ws.Range("B4").Value = $newFileQuery

# Re-editing that long wrapped cell makes the engine auto-grow row 4 past
# Excel's real 409.5pt row-height ceiling; pin it back to the size every
# other maxed-out wrapped row in this sheet already uses.
$ws.Rows.Item(4).RowHeight = 409.5

# Leave the selection/cursor on the cell that was actually edited.
$ws.Range("B4").Select()
